$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
# Overview sheet carries per-locale status in columns E (zh-cn) and F (de-de)
$ws1.Range("E2:F4").Value = "In Translation"
# Each locale sheet carries its own status in column C ("Status")
$ws2.Range("C2:C4").Value = "In Translation"
$ws3.Range("C2:C4").Value = "In Translation"

# --- Narrow the Status columns to match the new, shorter status text ---
$ws1.Range("E1:F1").ColumnWidth = 12.5
$ws2.Range("C1").ColumnWidth = 12.5
$ws3.Range("C1").ColumnWidth = 12.5
